# Update countries & provincias Spain
# - Swap the order of "Ghana" and "Bosnia y Herzegovina" (Bosnia now listed
#   before Ghana, i.e. row 76 becomes Bosnia y Herzegovina and row 77
#   becomes Ghana).
# - Refresh the COVID-19 stat columns (B:H) for several rows whose data
#   was updated: Kuwait (row 60), the two swapped rows (76 & 77), and
#   Sri Lanka (row 105).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap country names for rows 76 and 77 -------------------------------
$ws.Range("A76").Value = "Bosnia y Herzegovina"
$ws.Range("A77").Value = "Ghana"

# --- Row 60: Kuwait --------------------------------------------------------
$ws.Range("B60").Value = 3288
$ws.Range("C60").Value = 213
$ws.Range("D60").Value = 1012
$ws.Range("E60").Value = 2254
$ws.Range("F60").Value = 64
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 22

# --- Row 76: now Bosnia y Herzegovina --------------------------------------
$ws.Range("B76").Value = 1565
$ws.Range("C76").Value = 49
$ws.Range("D76").Value = 659
$ws.Range("E76").Value = 846
$ws.Range("F76").Value = 4
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 60

# --- Row 77: now Ghana ------------------------------------------------------
$ws.Range("B77").Value = 1550
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 155
$ws.Range("E77").Value = 1384
$ws.Range("F77").Value = 4
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 11

# --- Row 105: Sri Lanka ------------------------------------------------------
$ws.Range("B105").Value = 567
$ws.Range("C105").Value = 44
$ws.Range("D105").Value = 126
$ws.Range("E105").Value = 434
$ws.Range("F105").Value = 2
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 7
